$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PSA IBP Adapter Cables")

$ws.Range("A4").Value = 2
$ws.Range("E4").Value = "Fogg System 0395-2319:L12"
$ws.Range("F4").Value = "Edwards TruWave to unterminated cable"
$ws.Range("G4").Value = "Can be replaced with any cable that has the gray Truwave connector on one end"

$ws.Activate()
$ws.Range("A5:G6").Select()
